# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
# Source diff adds row 39: 2025-09-25 / 15:22:43 / "1.00 EUR = 1,621.3766"
# and bumps the sheet dimension from A1:C38 to A1:C39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 39

# Use a leading apostrophe so Excel stores the date-/time-looking strings
# as plain text (matching the existing inline-string cells) instead of
# auto-converting them into date/time serial numbers.
$ws.Cells.Item($newRow, 1).Value = "'2025-09-25"
$ws.Cells.Item($newRow, 2).Value = "'15:22:43"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,621.3766"

# Reset the cell style so the forced text-entry (quote prefix) doesn't leave
# behind a distinct number-format/style compared to the rest of the sheet.
$ws.Cells.Item($newRow, 1).Style = "Normal"
$ws.Cells.Item($newRow, 2).Style = "Normal"
